# Apply the row-data permutation among rows 24-28 on the "Artfynd" sheet.
# The location-fixed columns (C, J, K, N, P, S, T, U, V, W, Y, AA, AD, AE,
# AF, AG, AH, AT, AW, AX, AY, ...) stay where they are; only the
# observation-identifying columns (A, B, D, E, F, G, H, I, Q, R) are
# shuffled between the five rows as follows (target row <- source row):
#   24 <- 25
#   25 <- 24
#   26 <- 27
#   27 <- 28
#   28 <- 26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "moving" columns for rows 24-28 before touching anything,
# so writes to one row never clobber data still needed for another row.
$cols = @("A","B","D","E","F","G","H","I","Q","R")
$rows = 24,25,26,27,28

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# target row -> source row
$mapping = @{ 24 = 25; 25 = 24; 26 = 27; 27 = 28; 28 = 26 }

foreach ($target in $rows) {
    $source = $mapping[$target]
    $src = $snapshot[$source]

    $ws.Range("A$target").Value = $src["A"]
    $ws.Range("B$target").Value = $src["B"]
    $ws.Range("D$target").Value = $src["D"]
    $ws.Range("E$target").Value = $src["E"]
    $ws.Range("F$target").Value = $src["F"]
    $ws.Range("G$target").Value = $src["G"]
    $ws.Range("H$target").Value = $src["H"]

    $iVal = $src["I"]
    if ($null -eq $iVal -or $iVal -eq "") {
        $ws.Range("I$target").ClearContents()
    } else {
        # "Antal" is stored as text even though it looks numeric - use a
        # leading apostrophe so Excel keeps it as text, matching the
        # original inlineStr representation.
        $ws.Range("I$target").Value = "'" + $iVal
    }

    $ws.Range("Q$target").Value = $src["Q"]
    $ws.Range("R$target").Value = $src["R"]
}
